# Generate Report for handoff
#
# The localization pipeline ran again: the old source file
# (95da1b9c-6255-441c-8bbe-4ca5f9cfd697.md) was replaced by a new one
# (c414689a-9b01-4383-8d03-75adc4629023.md), its handoff package hash
# changed (72ae689c... -> 0b290065...) and new handoff timestamps were
# recorded. The previous "Handoff transform failed" row
# (f8b53712-a37b-43cb-a21e-f5e383ace1f2.md) is gone from the report,
# and the ".localization-config" row moves up to take its place.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($sheet, $row, $col) {
    $found = $true
    while ($found) {
        $found = $false
        foreach ($h in $sheet.Hyperlinks) {
            if ($h.Range.Row -eq $row -and $h.Range.Column -eq $col) {
                $h.Delete()
                $found = $true
                break
            }
        }
    }
}

function Set-HyperlinkDisplay($sheet, $row, $col, $text) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Row -eq $row -and $h.Range.Column -eq $col) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop the old row 4 (it held the ".localization-config" entry, which
# moves up into row 3 below instead).
$ws1.Rows(4).Delete()

# Row 3 used to describe f8b53712-...md / "Handoff transform failed".
# That file no longer appears in the report; the row now carries the
# ".localization-config" / "Not to be localized" entry instead (its
# hyperlink is kept in place and simply relabeled below).
# Deleting row 4 leaves its hyperlink (originally ".localization-config")
# dangling at the now-nonexistent row 4; drop it.
Remove-HyperlinkAt $ws1 4 1
$ws1.Range("A3").Value = ".localization-config"
$ws1.Range("B3").Value = "Not to be localized"
$ws1.Range("C3").Value = "Not to be localized"

# Row 2 now reports on the new source file.
$ws1.Range("A2").Value = "c414689a-9b01-4383-8d03-75adc4629023.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

Set-HyperlinkDisplay $ws1 2 1 "c414689a-9b01-4383-8d03-75adc4629023.md"
Set-HyperlinkDisplay $ws1 3 1 ".localization-config"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows(4).Delete()

Remove-HyperlinkAt $ws2 4 1
$ws2.Range("A3").Value = ".localization-config"
$ws2.Range("B3").Value = "Not to be localized"

$ws2.Range("A2").Value = "c414689a-9b01-4383-8d03-75adc4629023.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "c414689a-9b01-4383-8d03-75adc4629023.0b290065dd64f57d9d7bbaf097bb0f942c8363f0.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-08 18:33:16"

Set-HyperlinkDisplay $ws2 2 1 "c414689a-9b01-4383-8d03-75adc4629023.md"
Set-HyperlinkDisplay $ws2 2 3 "c414689a-9b01-4383-8d03-75adc4629023.0b290065dd64f57d9d7bbaf097bb0f942c8363f0.zh-cn.xlf"
Set-HyperlinkDisplay $ws2 3 1 ".localization-config"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows(4).Delete()

Remove-HyperlinkAt $ws3 4 1
$ws3.Range("A3").Value = ".localization-config"
$ws3.Range("B3").Value = "Not to be localized"

$ws3.Range("A2").Value = "c414689a-9b01-4383-8d03-75adc4629023.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "c414689a-9b01-4383-8d03-75adc4629023.0b290065dd64f57d9d7bbaf097bb0f942c8363f0.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-08 18:33:25"

Set-HyperlinkDisplay $ws3 2 1 "c414689a-9b01-4383-8d03-75adc4629023.md"
Set-HyperlinkDisplay $ws3 2 3 "c414689a-9b01-4383-8d03-75adc4629023.0b290065dd64f57d9d7bbaf097bb0f942c8363f0.de-de.xlf"
Set-HyperlinkDisplay $ws3 3 1 ".localization-config"
